# T460 update: refresh trip-sheet data (new km readings, trips, and locations)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Km initiali (start of period odometer reading)
$ws.Cells.Item(12, 2).Value = 229926

# Daily rows: r => (Km_parcursi, Locul deplasarii, Observatii utilizator)
$ws.Cells.Item(15, 2).Value = 421
$ws.Cells.Item(15, 3).Value = "Cluj-Satu-Mare"
$ws.Cells.Item(15, 4).Value = "Interes Serviciu"

$ws.Cells.Item(16, 2).Value = 30
$ws.Cells.Item(16, 3).Value = "Acasa-Birou"
$ws.Cells.Item(16, 4).Value = " "

$ws.Cells.Item(17, 2).Value = 30
$ws.Cells.Item(17, 3).Value = "Acasa-Birou"
$ws.Cells.Item(17, 4).Value = " "

$ws.Cells.Item(19, 2).Value = 30
$ws.Cells.Item(19, 3).Value = "Acasa-Birou"
$ws.Cells.Item(19, 4).Value = " "

$ws.Cells.Item(22, 2).Value = 421
$ws.Cells.Item(22, 3).Value = "Cluj-Satu-Mare"
$ws.Cells.Item(22, 4).Value = "Interes Serviciu"

$ws.Cells.Item(23, 2).Value = 101
$ws.Cells.Item(23, 3).Value = "Cluj-Dej"
$ws.Cells.Item(23, 4).Value = "Interes Serviciu"

$ws.Cells.Item(24, 2).Value = 156
$ws.Cells.Item(24, 3).Value = "Cluj-Zalau"
$ws.Cells.Item(24, 4).Value = "Interes Serviciu"

$ws.Cells.Item(25, 2).Value = 92
$ws.Cells.Item(25, 3).Value = "Cluj-Bontida"
$ws.Cells.Item(25, 4).Value = "Interes Serviciu"

$ws.Cells.Item(26, 2).Value = 30
$ws.Cells.Item(26, 3).Value = "Acasa-Birou"
$ws.Cells.Item(26, 4).Value = " "

$ws.Cells.Item(29, 2).Value = 30
$ws.Cells.Item(29, 3).Value = "Acasa-Birou"
$ws.Cells.Item(29, 4).Value = " "

$ws.Cells.Item(30, 2).Value = 421
$ws.Cells.Item(30, 3).Value = "Cluj-Satu-Mare"
$ws.Cells.Item(30, 4).Value = "Interes Serviciu"

$ws.Cells.Item(31, 2).Value = 92
$ws.Cells.Item(31, 3).Value = "Cluj-Bontida"
$ws.Cells.Item(31, 4).Value = "Interes Serviciu"

$ws.Cells.Item(32, 2).Value = 421
$ws.Cells.Item(32, 3).Value = "Cluj-Satu-Mare"
$ws.Cells.Item(32, 4).Value = "Interes Serviciu"

$ws.Cells.Item(33, 2).Value = 30
$ws.Cells.Item(33, 3).Value = "Acasa-Birou"
$ws.Cells.Item(33, 4).Value = " "

$ws.Cells.Item(36, 2).Value = 30
$ws.Cells.Item(36, 3).Value = "Acasa-Birou"
$ws.Cells.Item(36, 4).Value = " "

$ws.Cells.Item(37, 2).Value = 121
$ws.Cells.Item(37, 3).Value = "Cluj-Turda"
$ws.Cells.Item(37, 4).Value = "Interes Serviciu"

$ws.Cells.Item(38, 2).Value = 30
$ws.Cells.Item(38, 3).Value = "Acasa-Birou"
$ws.Cells.Item(38, 4).Value = " "

$ws.Cells.Item(39, 2).Value = 121
$ws.Cells.Item(39, 3).Value = "Cluj-Turda"
$ws.Cells.Item(39, 4).Value = "Interes Serviciu"

$ws.Cells.Item(40, 2).Value = 30
$ws.Cells.Item(40, 3).Value = "Acasa-Birou"
$ws.Cells.Item(40, 4).Value = " "

$ws.Cells.Item(43, 2).Value = 30
$ws.Cells.Item(43, 3).Value = "Acasa-Birou"
$ws.Cells.Item(43, 4).Value = " "

$ws.Cells.Item(44, 2).Value = 92
$ws.Cells.Item(44, 3).Value = "Cluj-Bontida"
$ws.Cells.Item(44, 4).Value = "Interes Serviciu"

# Totals
$ws.Cells.Item(45, 2).Value = 2759
$ws.Cells.Item(46, 2).Value = 232685
